$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

Set-TextValue $ws.Range('E3') '2'
Set-TextValue $ws.Range('F3') '203 / 4,901'
Set-TextValue $ws.Range('G3') '6,503,000'
Set-TextValue $ws.Range('H3') '8,670,609'
Set-TextValue $ws.Range('I3') '30,540'
Set-TextValue $ws.Range('K3') '265,592'
Set-TextValue $ws.Range('L3') '07-24 15:27'
Set-TextValue $ws.Range('F4') '63 / 5,093'
Set-TextValue $ws.Range('G4') '1,907,128'
Set-TextValue $ws.Range('H4') '129,884'
Set-TextValue $ws.Range('I4') '2,923'
Set-TextValue $ws.Range('J4') '2,139'
Set-TextValue $ws.Range('K4') '174,417'
Set-TextValue $ws.Range('L4') '07-24 15:27'
Set-TextValue $ws.Range('E5') '1'
Set-TextValue $ws.Range('F5') '195 / 6,788'
Set-TextValue $ws.Range('G5') '1,000,898'
Set-TextValue $ws.Range('H5') '1,443,566'
Set-TextValue $ws.Range('I5') '8,032'
Set-TextValue $ws.Range('K5') '333,297'
Set-TextValue $ws.Range('L5') '07-24 15:28'
Set-TextValue $ws.Range('E6') '1,314'
Set-TextValue $ws.Range('G6') '5,788,217'
Set-TextValue $ws.Range('H6') '1,173,326'
Set-TextValue $ws.Range('I6') '1,442'
Set-TextValue $ws.Range('J6') '4,016'
Set-TextValue $ws.Range('K6') '148,689'
Set-TextValue $ws.Range('L6') '07-24 15:18'
Set-TextValue $ws.Range('G7') '175,229'
Set-TextValue $ws.Range('H7') '57,499'
Set-TextValue $ws.Range('I7') '1,025'
Set-TextValue $ws.Range('J7') '15,208'
Set-TextValue $ws.Range('K7') '179,215'
Set-TextValue $ws.Range('L7') '07-24 15:21'
Set-TextValue $ws.Range('E8') '593'
Set-TextValue $ws.Range('F8') '77 / 5,188'
Set-TextValue $ws.Range('G8') '921,233'
Set-TextValue $ws.Range('H8') '1,836,341'
Set-TextValue $ws.Range('I8') '1,435'
Set-TextValue $ws.Range('J8') '119'
Set-TextValue $ws.Range('K8') '25,476'
Set-TextValue $ws.Range('L8') '07-24 15:22'
Set-TextValue $ws.Range('E9') '96'
Set-TextValue $ws.Range('G9') '336,331'
Set-TextValue $ws.Range('H9') '245,940'
Set-TextValue $ws.Range('I9') '427'
Set-TextValue $ws.Range('K9') '21,457'
Set-TextValue $ws.Range('L9') '07-24 15:28'
Set-TextValue $ws.Range('E10') '1'
Set-TextValue $ws.Range('F10') '24 / 2,500'
Set-TextValue $ws.Range('G10') '7,820,442'
Set-TextValue $ws.Range('H10') '7,066,519'
Set-TextValue $ws.Range('I10') '333,258'
Set-TextValue $ws.Range('J10') '299'
Set-TextValue $ws.Range('K10') '204,255'
Set-TextValue $ws.Range('L10') '07-24 15:27'
Set-TextValue $ws.Range('F11') '3 / 1,621'
Set-TextValue $ws.Range('G11') '9,280,123'
Set-TextValue $ws.Range('H11') '1,511,124'
Set-TextValue $ws.Range('I11') '5,317'
Set-TextValue $ws.Range('K11') '73,477'
Set-TextValue $ws.Range('L11') '07-24 15:27'
Set-TextValue $ws.Range('F12') '15 / 1,378'
Set-TextValue $ws.Range('G12') '1,004,118'
Set-TextValue $ws.Range('H12') '540,008'
Set-TextValue $ws.Range('I12') '591'
Set-TextValue $ws.Range('K12') '13,743'
Set-TextValue $ws.Range('L12') '07-24 15:27'
Set-TextValue $ws.Range('H13') '521,343'
Set-TextValue $ws.Range('I13') '16'
Set-TextValue $ws.Range('K13') '39,056'
Set-TextValue $ws.Range('L13') '07-24 15:27'
Set-TextValue $ws.Range('G14') '179,669'
Set-TextValue $ws.Range('H14') '36,396'
Set-TextValue $ws.Range('I14') '30'
Set-TextValue $ws.Range('K14') '6,498'
Set-TextValue $ws.Range('L14') '07-24 15:27'
Set-TextValue $ws.Range('A15') 'rossmann-de'
Set-TextValue $ws.Range('F15') '407'
Set-TextValue $ws.Range('G15') '72,156'
Set-TextValue $ws.Range('H15') '441,075'
Set-TextValue $ws.Range('I15') '174'
Set-TextValue $ws.Range('J15') '0'
Set-TextValue $ws.Range('K15') '10,210'
Set-TextValue $ws.Range('L15') '07-24 15:27'
Set-TextValue $ws.Range('A16') 'rossmann-pl'
Set-TextValue $ws.Range('F16') '6 / 483'
Set-TextValue $ws.Range('G16') '56,039'
Set-TextValue $ws.Range('H16') '18,884'
Set-TextValue $ws.Range('I16') '11'
Set-TextValue $ws.Range('J16') '24'
Set-TextValue $ws.Range('K16') '894'
Set-TextValue $ws.Range('L16') '07-24 15:27'
Set-TextValue $ws.Range('E17') '267'
Set-TextValue $ws.Range('G17') '77,168'
Set-TextValue $ws.Range('H17') '41,444'
Set-TextValue $ws.Range('I17') '1'
Set-TextValue $ws.Range('L17') '07-24 15:27'
Set-TextValue $ws.Range('E18') '35'
Set-TextValue $ws.Range('G18') '2,439,136'
Set-TextValue $ws.Range('H18') '518,300'
Set-TextValue $ws.Range('I18') '1,640'
Set-TextValue $ws.Range('K18') '6,856'
Set-TextValue $ws.Range('L18') '07-24 15:27'
Set-TextValue $ws.Range('G19') '137,454'
Set-TextValue $ws.Range('H19') '63,637'
Set-TextValue $ws.Range('I19') '6'
Set-TextValue $ws.Range('L19') '07-24 15:26'
Set-TextValue $ws.Range('E20') '30'
Set-TextValue $ws.Range('F20') '18 / 219'
Set-TextValue $ws.Range('G20') '429,937'
Set-TextValue $ws.Range('H20') '135,321'
Set-TextValue $ws.Range('I20') '31'
Set-TextValue $ws.Range('K20') '6,489'
Set-TextValue $ws.Range('L20') '07-24 15:27'
Set-TextValue $ws.Range('F21') '2 / 490'
Set-TextValue $ws.Range('G21') '208,169'
Set-TextValue $ws.Range('H21') '11,599'
Set-TextValue $ws.Range('I21') '563'
Set-TextValue $ws.Range('J21') '63'
Set-TextValue $ws.Range('K21') '22,579'
Set-TextValue $ws.Range('L21') '07-24 15:27'
Set-TextValue $ws.Range('G22') '324,985'
Set-TextValue $ws.Range('H22') '163,932'
Set-TextValue $ws.Range('I22') '142'
Set-TextValue $ws.Range('K22') '1,385'
Set-TextValue $ws.Range('L22') '07-24 15:27'
Set-TextValue $ws.Range('H23') '37,702'
Set-TextValue $ws.Range('I23') '56'
Set-TextValue $ws.Range('K23') '2,959'
Set-TextValue $ws.Range('L23') '07-24 15:26'
Set-TextValue $ws.Range('F24') '9 / 324'
Set-TextValue $ws.Range('H24') '54,612'
Set-TextValue $ws.Range('I24') '220'
Set-TextValue $ws.Range('K24') '17,454'
Set-TextValue $ws.Range('L24') '07-24 15:27'
Set-TextValue $ws.Range('G25') '1,774,562'
Set-TextValue $ws.Range('H25') '1,292'
Set-TextValue $ws.Range('I25') '9'
Set-TextValue $ws.Range('L25') '07-24 15:26'
Set-TextValue $ws.Range('G26') '147,548'
Set-TextValue $ws.Range('H26') '145,012'
Set-TextValue $ws.Range('I26') '989'
Set-TextValue $ws.Range('K26') '5,259'
Set-TextValue $ws.Range('L26') '07-24 15:26'
Set-TextValue $ws.Range('G27') '671'
Set-TextValue $ws.Range('H27') '1,106'
Set-TextValue $ws.Range('I27') '5'
Set-TextValue $ws.Range('L27') '07-24 15:26'
Set-TextValue $ws.Range('G28') '5,925'
Set-TextValue $ws.Range('H28') '4,920'
Set-TextValue $ws.Range('I28') '0'
Set-TextValue $ws.Range('K28') '1,181'
Set-TextValue $ws.Range('L28') '07-24 15:26'
Set-TextValue $ws.Range('H29') '4,440'
Set-TextValue $ws.Range('I29') '1'
Set-TextValue $ws.Range('K29') '785'
Set-TextValue $ws.Range('L29') '07-24 15:26'
Set-TextValue $ws.Range('G30') '202,136'
Set-TextValue $ws.Range('H30') '900'
Set-TextValue $ws.Range('L30') '07-24 15:26'
Set-TextValue $ws.Range('E31') '0'
Set-TextValue $ws.Range('G31') '2,764'
Set-TextValue $ws.Range('H31') '80,329'
Set-TextValue $ws.Range('I31') '4'
Set-TextValue $ws.Range('K31') '203'
Set-TextValue $ws.Range('L31') '07-24 15:26'
Set-TextValue $ws.Range('H32') '12,637'
Set-TextValue $ws.Range('K32') '135'
Set-TextValue $ws.Range('L32') '07-24 15:26'
Set-TextValue $ws.Range('H33') '3,717'
Set-TextValue $ws.Range('I33') '1'
Set-TextValue $ws.Range('L33') '07-24 15:26'
Set-TextValue $ws.Range('H34') '706'
Set-TextValue $ws.Range('L34') '07-24 15:26'
Set-TextValue $ws.Range('G35') '5,776'
Set-TextValue $ws.Range('H35') '3,650'
Set-TextValue $ws.Range('L35') '07-24 15:26'
Set-TextValue $ws.Range('L36') '07-24 15:26'
Set-TextValue $ws.Range('L38') '07-24 15:21'
Set-TextValue $ws.Range('L39') '07-24 15:26'
